# Remove the subtitle placeholder ("副標題 2" / type=subTitle idx=1, which
# holds the "Peter H. Chen" byline) from the section-divider slides 3, 6 and
# 9. Slide 1 (the top-level "4 Architecture" divider) keeps its copy.
#
# NOTE: deleting a slide-level placeholder shape only clears/resets it back
# to an empty layout-inherited placeholder on the first call to Delete() -
# the Shapes collection still reports the same count, just with a new Id
# and an auto-generated English name ("Subtitle N"). A second Delete() call
# on that now-empty inherited placeholder actually removes the <p:sp> from
# the slide. We loop (matching on the placeholder's format Type, which is
# stable across the reset) until it is well and truly gone.

$p = $ppt.ActivePresentation
$slideIndexes = @(3, 6, 9)
$ppPlaceholderSubtitle = 4

foreach ($si in $slideIndexes) {
    $s = $p.Slides.Item($si)

    $stillThere = $true
    while ($stillThere) {
        $stillThere = $false
        for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
            $sh = $s.Shapes.Item($i)
            if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderSubtitle) {
                $sh.Delete()
                $stillThere = $true
            }
        }
    }
}
